$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C2").Value = "ola"
$ws.Range("C3").Value = "ola"
$ws.Range("C4").Value = "ola"
